$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Periodo Mora" column (E) values: shift periods forward by one
# 2507/2506/2505 -> 2506/2507/2508
$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2507"
$ws.Range("E18").Value = "2508"

# Update "Salario Basico" column (G) values for the three rows
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
